$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Total money for winners" column (D) values from 100000 to 120000
# for data rows 2 through 13.
$ws.Range("D2:D13").Value = 120000

# Update the active selection to match the saved view state.
$ws.Range("F12").Select()
